$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (new Police Commissioner, volume/date range)
$ws.Range("M6").Value = "Edward A. Caban"
$ws.Range("A8").Value = "Volume 30   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  7/3/2023  Through  7/9/2023"

# Weekly crime statistics table updates
# Row 14
$ws.Range("C14").NumberFormat = '#,##0'
$ws.Range("C14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 5
$ws.Range("K14").Value = -16.666666666666
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -16.666666666666
$ws.Range("N14").Value = -82.758620689655
# Row 15
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = -23.076923076923
$ws.Range("M15").Value = -23.076923076923
$ws.Range("N15").Value = -79.166666666666
# Row 16
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -27.272727272727
$ws.Range("F16").Value = 41
$ws.Range("G16").Value = 44
$ws.Range("H16").Value = -6.818181818181
$ws.Range("I16").Value = 172
$ws.Range("J16").Value = 214
$ws.Range("K16").Value = -19.626168224299
$ws.Range("L16").Value = -10.416666666666
$ws.Range("M16").Value = -18.867924528301
$ws.Range("N16").Value = -84.724689165186
# Row 17
$ws.Range("C17").Value = 26
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = 73.333333333333
$ws.Range("F17").Value = 78
$ws.Range("G17").Value = 68
$ws.Range("H17").Value = 14.705882352941
$ws.Range("I17").Value = 394
$ws.Range("J17").Value = 372
$ws.Range("K17").Value = 5.913978494623
$ws.Range("L17").Value = 26.282051282051
$ws.Range("M17").Value = 79.908675799086
$ws.Range("N17").Value = -35.830618892508
# Row 18
$ws.Range("F18").Value = 15
$ws.Range("H18").Value = -31.818181818181
$ws.Range("I18").Value = 82
$ws.Range("J18").Value = 111
$ws.Range("K18").Value = -26.126126126126
$ws.Range("L18").Value = -17.171717171717
$ws.Range("M18").Value = -21.904761904761
$ws.Range("N18").Value = -92.949269131556
# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 44.444444444444
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 13.043478260869
$ws.Range("I19").Value = 282
$ws.Range("J19").Value = 323
$ws.Range("K19").Value = -12.693498452012
$ws.Range("L19").Value = -7.843137254901
$ws.Range("M19").Value = 83.116883116883
$ws.Range("N19").Value = -42.448979591836
# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -35.294117647058
$ws.Range("I20").Value = 146
$ws.Range("J20").Value = 118
$ws.Range("K20").Value = 23.728813559322
$ws.Range("L20").Value = 111.594202898551
$ws.Range("M20").Value = 124.615384615385
$ws.Range("N20").Value = -71.705426356589
# Row 21
$ws.Range("C21").Value = 54
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = 20
$ws.Range("F21").Value = 198
$ws.Range("G21").Value = 200
$ws.Range("H21").Value = -1
$ws.Range("I21").Value = 1091
$ws.Range("J21").Value = 1157
$ws.Range("K21").Value = -5.704407951598
$ws.Range("L21").Value = 9.758551307847
$ws.Range("M21").Value = 40.956072351421
$ws.Range("N21").Value = -72.629202207727
# Row 22
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 1
$ws.Range("I22").Value = 11
$ws.Range("K22").Value = -15.384615384615
$ws.Range("L22").Value = -50
$ws.Range("M22").Value = -21.428571428571
# Row 23
$ws.Range("D23").Value = 2
$ws.Range("G23").Value = 6
$ws.Range("J23").Value = 18
$ws.Range("K23").Value = -11.111111111111
$ws.Range("M23").Value = 6.666666666666
# Row 24
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 21.739130434782
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = -15.384615384615
$ws.Range("I24").Value = 542
$ws.Range("J24").Value = 727
$ws.Range("K24").Value = -25.447042640990
$ws.Range("L24").Value = 9.716599190283
$ws.Range("M24").Value = 11.065573770491
# Row 25
$ws.Range("C25").Value = 38
$ws.Range("D25").Value = 39
$ws.Range("E25").Value = -2.564102564102
$ws.Range("F25").Value = 126
$ws.Range("G25").Value = 121
$ws.Range("H25").Value = 4.132231404958
$ws.Range("I25").Value = 672
$ws.Range("J25").Value = 554
$ws.Range("K25").Value = 21.299638989169
$ws.Range("L25").Value = 51.011235955056
$ws.Range("M25").Value = -1.321585903083
# Row 26
$ws.Range("D26").Value = 2
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -75
$ws.Range("J26").Value = 25
$ws.Range("K26").Value = 20
# Row 27
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 38
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = -24
$ws.Range("L27").Value = -13.636363636363
# Row 28
$ws.Range("C28").Value = 4
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 275
$ws.Range("I28").Value = 29
$ws.Range("K28").Value = 26.086956521739
$ws.Range("L28").Value = 31.818181818181
$ws.Range("M28").Value = 16
$ws.Range("N28").Value = -68.131868131868
# Row 29
$ws.Range("C29").Value = 3
$ws.Range("F29").Value = 11
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 175
$ws.Range("I29").Value = 22
$ws.Range("K29").Value = 4.761904761904
$ws.Range("L29").Value = 4.761904761904
$ws.Range("M29").Value = 4.761904761904
$ws.Range("N29").Value = -74.117647058823
